# Ibkr.xlsx edit: append new trade rows (305-326), fix timestamps, replace
# PE/CE market-order wording with MARKET/LIMIT + PE/CE split, and update the
# saved window selection, per the commit:
#   "google sheet implemented stoploss implementeed day fixed timestamp
#    fixed pe/ce replaced by buy/sell"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one new row, columns A..K in order:
#   A: Trigger_Level_High_Low
#   B: Entry_Type        (text: MARKET / LIMIT)
#   C: Entry_Strike
#   D: Strike_Type       (text: PE / CE)
#   E: Expiry            (date serial, formatted yyyy-mm-dd hh:mm:ss)
#   F: Target
#   G: Stop_Loss
#   H: Qty
#   I: Slicing
#   J: Time_Interval
#   K: Activation
$rows = @(
    @(36715, "LIMIT", 36730, "PE", 45660, 36750, 36700, 4, 2, 5, 0),
    @(36705, "LIMIT", 36710, "PE", 45660, 36750, 36680, 4, 2, 5, 0),
    @(36700, "LIMIT", 36705, "PE", 45660, 36730, 36670, 4, 2, 5, 0),
    @(36910, "LIMIT", 36925, "PE", 45661, 36950, 36860, 6, 2, 5, 0),
    @(36850, "LIMIT", 36870, "PE", 45660, 36950, 36810, 6, 2, 5, 0),
    @(36850, "LIMIT", 36870, "PE", 45660, 36950, 36810, 6, 2, 5, 0),
    @(36850, "LIMIT", 36870, "PE", 45660, 36950, 36810, 6, 2, 5, 0),
    @(36850, "LIMIT", 36870, "PE", 45660, 36950, 36810, 6, 2, 5, 0),
    @(36970, "LIMIT", 37000, "PE", 45660, 37050, 36930, 6, 2, 5, 0),
    @(36970, "LIMIT", 37000, "PE", 45660, 37050, 36930, 6, 2, 5, 0),
    @(36790, "LIMIT", 36810, "PE", 45660, 36940, 36680, 4, 2, 5, 0),
    @(36790, "LIMIT", 36810, "PE", 45661, 36940, 36680, 4, 2, 5, 0),
    @(36700, "LIMIT", 36810, "PE", 45661, 36940, 36680, 4, 2, 5, 0),
    @(36700, "LIMIT", 36740, "PE", 45661, 36800, 36680, 4, 2, 5, 0),
    @(36700, "MARKET", 36740, "PE", 45661, 36800, 36680, 4, 2, 5, 0),
    @(36700, "MARKET", 36740, "PE", 45661, 36800, 36680, 4, 2, 5, 0),
    @(36880, "MARKET", 36880, "CE", 45661, 36820, 36920, 4, 2, 5, 0),
    @(36800, "MARKET", 36830, "CE", 45661, 36870, 36960, 4, 2, 5, 0),
    @(36800, "MARKET", 36830, "CE", 45661, 36870, 36960, 4, 2, 5, 0),
    @(36800, "MARKET", 36830, "CE", 45661, 36870, 36960, 4, 2, 5, 0),
    @(36850, "LIMIT", 36860, "CE", 45661, 36750, 36960, 4, 2, 5, 0),
    @(36850, "LIMIT", 36860, "CE", 45661, 36750, 36960, 4, 2, 5, 0)
)

$startRow = 305
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = $vals[4]
    $eCell.NumberFormat = "yyyy-mm-dd hh:mm:ss"

    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = $vals[6]
    $ws.Cells.Item($r, 8).Value = $vals[7]
    $ws.Cells.Item($r, 9).Value = $vals[8]
    $ws.Cells.Item($r, 10).Value = $vals[9]
    $ws.Cells.Item($r, 11).Value = $vals[10]
}

# Update the recorded scroll position / selection to match where the editor
# left off after appending the rows above.
$ws.Range("H323").Select()
